$d = $word.ActiveDocument
$W = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------------
# 1) Paragraph 5 ("Add "There are no speaker notes..."): collapse the
#    five separate runs into a single run with identical text.
# ---------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5)
$r5 = $p5.Range
$r5.End = $r5.End - 1
$r5.Text = $r5.Text + " "
$r5b = $d.Paragraphs.Item(5).Range
$r5b.End = $r5b.End - 1
$r5b.Text = 'Add "There are no speaker notes on this slide." for any blank slides'

# ---------------------------------------------------------------------
# 2) Remove the stray _GoBack bookmark that currently sits on the
#    "Convert superscripts..." paragraph - it is being relocated.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 3) Paragraph 16 ("Find what: Slide [0-9]{1,}"): re-case "Slide" to
#    "SLIDE" (split across three runs) and re-add the _GoBack bookmark
#    at the end of the paragraph.
# ---------------------------------------------------------------------
$p16 = $d.Paragraphs.Item(16)
$r16 = $p16.Range
$r16.End = $r16.End - 1
$r16.Text = ""
$xml16 = '<w:p xmlns:w="' + $W + '">' `
    + '<w:r><w:t xml:space="preserve">Find what: </w:t></w:r>' `
    + '<w:r><w:t>SLIDE</w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve"> [0-9]{1,}</w:t></w:r>' `
    + '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' `
    + '<w:bookmarkEnd w:id="0"/>' `
    + '</w:p>'
[void]$r16.InsertXML($xml16)

# ---------------------------------------------------------------------
# 4) Paragraph 18 ("Replace with: ',^p'"): re-split the runs and wrap
#    the comma/caret in gramStart/gramEnd proofErr markers.
# ---------------------------------------------------------------------
$p18 = $d.Paragraphs.Item(18)
$r18 = $p18.Range
$r18.End = $r18.End - 1
$r18.Text = ""
$xml18 = '<w:p xmlns:w="' + $W + '">' `
    + '<w:r><w:t xml:space="preserve">Replace with: </w:t></w:r>' `
    + "<w:r><w:t>'</w:t></w:r>" `
    + '<w:proofErr w:type="gramStart"/>' `
    + '<w:r><w:t>,</w:t></w:r>' `
    + '<w:r><w:t>^</w:t></w:r>' `
    + '<w:proofErr w:type="gramEnd"/>' `
    + '<w:r><w:t>p</w:t></w:r>' `
    + "<w:r><w:t>'</w:t></w:r>" `
    + '</w:p>'
[void]$r18.InsertXML($xml18)
